# Add a "330 Ohm Resistor" line item to the Arduino Nano BOM.
# This mirrors a user inserting a new row above the "9V Battery" row
# (row 10) and filling in Part / Quantity / Unit Price / Subtotal,
# matching the rest of the table's layout & formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10, pushing everything else down.
$ws.Rows.Item(10).Insert()

# Part name (new shared string)
$ws.Cells.Item(10, 1).Value2() = "330 Ohm Resistor"

# Quantity
$ws.Cells.Item(10, 2).Value2() = 1

# Unit Price
$ws.Cells.Item(10, 3).Value2() = 0.04

# Subtotal = Quantity * Unit Price (same pattern as the rest of column D)
$ws.Cells.Item(10, 4).Formula() = "=B10*C10"

Write-Output "Added 330 Ohm Resistor row to BOM"
